# Apply the "DaySale" refresh: 6 new out-of-stock products are inserted
# (alphabetically) into the shortage list, the existing rows shift their
# quantities accordingly, the grand total and the footer timestamp update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room: insert 6 blank rows right before the current Total row
#    (row 25), pushing Total -> row 31 and the footer -> row 32.
# ---------------------------------------------------------------------
$ws.Rows("25:30").Insert()

# Clone formatting/merges for the 6 new data rows from the last existing
# data block (rows 19:24 -> 25:30) so borders/merges/number formats match.
$ws.Range("A19:Q24").Copy($ws.Range("A25"))

# ---------------------------------------------------------------------
# 2) Final alphabetical product list (m = 1..24) with all column values:
#    A = item no., C = name, H = stock, L = order limit,
#    N = price, P = sell price, Q = transactions
# ---------------------------------------------------------------------
$data = @(
  @{Row=7;  A=1;  C="ANTI-COX II 15MG/3ML 6 AMP";                       H="1:2";  L="1"; N="78.00";  P="12.4800"; Q="0:1"},
  @{Row=8;  A=2;  C="BLOKATENS 10/160MG 28 F.C.TABS.";                  H="0:1";  L="1"; N="160.00"; P="160.0000";Q="1:0"},
  @{Row=9;  A=3;  C="CALCITRON 30 CAPS.";                               H="0:2";  L="1"; N="132.00"; P="43.5600"; Q="0:1"},
  @{Row=10; A=4;  C="CATAFLAM 50 MG 20 SUGAR C.TABS";                   H="1:1";  L="1"; N="86.00";  P="43.0000"; Q="0:1"},
  @{Row=11; A=5;  C="CONTAFEVER N 200MG/5ML SUSP. 120ML";               H="14:0"; L="1"; N="33.00";  P="33.0000"; Q="1:0"},
  @{Row=12; A=6;  C="CONVENTIN 300MG 30 CAPS.";                         H="1:0";  L="1"; N="114.00"; P="37.6200"; Q="0:1"},
  @{Row=13; A=7;  C="ERASTAPEX PLUS 20MG/12.5MG 30 TAB";                H="1:2";  L="1"; N="78.00";  P="25.7400"; Q="0:1"},
  @{Row=14; A=8;  C="FLECTOR 50MG 30 CAPS";                             H="1:0";  L="1"; N="87.00";  P="28.7100"; Q="0:1"},
  @{Row=15; A=9;  C="FLUMOX 500MG 16 CAPS";                             H="1:0";  L="1"; N="71.00";  P="35.5000"; Q="0:1"},
  @{Row=16; A=10; C="FORBUDES 400/12MCG 60 INHALATION CAPS.+INHALER";   H="1:1";  L="1"; N="334.00"; P="53.4400"; Q="0:1"},
  @{Row=17; A=11; C="GLUCOVANCE 500/5MG 30 F.C.TAB.";                   H="1:1";  L="1"; N="74.00";  P="37.0000"; Q="0:1"},
  @{Row=18; A=12; C="HYACARENOL EYE DROPS 10 ML";                       H="1:0";  L="1"; N="56.00";  P="56.0000"; Q="1:0"},
  @{Row=19; A=13; C="IVERZINE 1% LOTION 60 ML";                         H="2:0";  L="1"; N="52.00";  P="52.0000"; Q="1:0"},
  @{Row=20; A=14; C="MORALACT  TAB";                                    H="0:1";  L="0"; N="90.00";  P="45.0000"; Q="0:1"},
  @{Row=21; A=15; C="PK-MERZ 100MG 30 F.C. TAB";                        H="0:0";  L="1"; N="81.00";  P="81.0000"; Q="1:0"},
  @{Row=22; A=16; C="RELAXON 30 CAP";                                   H="2:1";  L="1"; N="69.00";  P="22.7700"; Q="0:1"},
  @{Row=23; A=17; C="SELGON 20MG 20 TABS.";                             H="18:1"; L="1"; N="30.00";  P="15.0000"; Q="0:1"},
  @{Row=24; A=18; C="SPASMOFREE 5MG/2ML I.V./I.M. 3 AMP.";              H="1:2";  L="1"; N="54.00";  P="17.8200"; Q="0:1"},
  @{Row=25; A=19; C="STATURIC 40MG 30 F.C. TABS.";                      H="0:1";  L="1"; N="93.00";  P="93.0000"; Q="1:0"},
  @{Row=26; A=20; C="TAMSULIN 0.4MG 28 CAPS";                           H="2:0";  L="1"; N="124.00"; P="62.0000"; Q="0:1"},
  @{Row=27; A=21; C="VASTAREL MR 35MG 30 F.C.TAB.";                     H="0:1";  L="1"; N="175.00"; P="175.0000";Q="1:0"},
  @{Row=28; A=22; C="VOLTAREN 75MG/3ML 3 AMP.";                         H="8:3";  L="1"; N="51.00";  P="16.8300"; Q="0:1"},
  @{Row=29; A=23; C="ZURCAL 40MG 14 GASTRO RESISTANT TAB";              H="4:0";  L="1"; N="96.00";  P="192.0000";Q="2:0"},
  @{Row=30; A=24; C="سرنجات 3 سم";                                      H="0:0";  L="0"; N="2.00";   P="6.0000";  Q="3:0"}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.A
    $ws.Range("C$r").Value = $item.C
    $ws.Range("H$r").Value = $item.H
    $ws.Range("L$r").Value = $item.L
    $ws.Range("N$r").Value = $item.N
    $ws.Range("P$r").Value = $item.P
    $ws.Range("Q$r").Value = $item.Q
}

# ---------------------------------------------------------------------
# 3) Grand total (sum of the "sell price" column) now at P31
# ---------------------------------------------------------------------
$ws.Range("P31").Value = 1344.47

# ---------------------------------------------------------------------
# 4) Footer timestamp (now row 32) reflects the new export time
# ---------------------------------------------------------------------
$ws.Range("A32").Value = "Sunday, 20 July, 2025 1:18 PM"
